# Regenerate merged AHB files
# - rename the "_old" / "_new" header suffixes to "_FV2410" / "_FV2504"
# - turn the data range into a real Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row: "_old" -> "_FV2410", "_new" -> "_FV2504"
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn A1:U66 into an Excel Table (created as "Table1")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U66"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)

# 3) Freeze the header row (split below row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
